# issue with vsync, gpu not being maximised
#
# The previous benchmark note ("." placeholder) is replaced everywhere by a
# bold correction, which had also been accidentally pasted across a wide
# swath of empty cells (columns D and O:P:Q) while the user was re-running
# benchmarks. A brand-new block of trial columns (pairs of "Seconds"/"FPS")
# was also added to the right of the last benchmark table (row 80 onward)
# to capture additional GPU runs.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$warning = "ALL WRONG, RTX GPU WAS NOT BEING USED AND NOT BEING MAXED OUT"

# The old "." note (previously only in Q29) is replaced by the warning text,
# which also got smeared across columns O:Q for rows 3-75 ...
$ws.Range("O3:Q75").Value = $warning

# ... and down column D for rows 52-67.
$ws.Range("D52:D67").Value = $warning

# New trial-header pairs ("Seconds"/"FPS") added alongside the existing
# A80:B80 pair, one pair per extra GPU run.
$headerPairs = @(
    @("D","E"), @("G","H"), @("J","K"), @("M","N"), @("P","Q"),
    @("S","T"), @("V","W"), @("Y","Z"), @("AB","AC")
)
foreach ($pair in $headerPairs) {
    $ws.Range($pair[0] + "80").Value = "Seconds"
    $ws.Range($pair[1] + "80").Value = "FPS"
}

# Data for the first three new trials (D:E, G:H, J:K) that were actually run.
$trial1 = 318,327,336,333,337,334,335,338,335,336
$trial2 = 316,338,333,332,330,317,302,295,324,311
$trial3 = 248,248,245,247,250,246,254,250,247,252

for ($i = 0; $i -lt 10; $i++) {
    $r = 81 + $i
    $ws.Cells.Item($r, 4).Value = $i + 1        # D: trial index
    $ws.Cells.Item($r, 5).Value = $trial1[$i]   # E: FPS
    $ws.Cells.Item($r, 7).Value = $i + 1        # G: trial index
    $ws.Cells.Item($r, 8).Value = $trial2[$i]   # H: FPS
    $ws.Cells.Item($r, 10).Value = $i + 1       # J: trial index
    $ws.Cells.Item($r, 11).Value = $trial3[$i]  # K: FPS
}

# A couple of stray leftover FPS readings below the J:K trial.
$ws.Range("K91").Value = 245
$ws.Range("K92").Value = 253

# The remaining new trial columns only got their index numbers filled in
# (1-10) -- the GPU runs for these were never actually recorded.
$indexOnlyCols = @("M","P","S","V","Y","AB")
foreach ($col in $indexOnlyCols) {
    for ($i = 0; $i -lt 10; $i++) {
        $ws.Range($col + (81 + $i)).Value = $i + 1
    }
}

# The saved view scrolled back to the top and left the selection on L62.
$ws.Range("L62").Select()
